# Updates cryptos list values (price + 1h volume change) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.390.80"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "2.232.34"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'245.16"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("D6").Value = "'0.629"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("D7").Value = "'74.01"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -3.97%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("D10").Value = "'43.08"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("D11").Value = "'0.0966"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("D12").Value = "'7.10"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "'14.41"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("E15").Value = "  -0.98%  "
$ws.Range("D16").Value = "2.236.45"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").Value = "42.180.43"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "'0.0000112"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +13.37%  "
$ws.Range("D19").Value = "'6.15"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("D20").Value = "'72.02"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").Value = "'10.28"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +38.08%  "
$ws.Range("D22").Value = "'230.90"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("E23").Value = "  -5.50%  "
$ws.Range("D24").Value = "'11.70"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.17%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "'3.69"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.06%  "
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("D28").Value = "'2.21"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.97%  "
$ws.Range("D29").Value = "'166.59"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D30").Value = "'20.88"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("D31").Value = "'5.83"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +17.60%  "
$ws.Range("D32").Value = "'0.0806"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.61%  "
$ws.Range("E33").Value = "  -2.39%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "'29.64"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -9.12%  "
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "'0.125"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("D36").Value = "'4.42"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.04%  "
$ws.Range("D37").Value = "'0.0307"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.55%  "
$ws.Range("D38").Value = "'13.18"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -8.11%  "
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("E40").Value = "  -4.03%  "
$ws.Range("D41").Value = "'63.04"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.52%  "
$ws.Range("D42").Value = "'0.200"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").Value = "'8.81"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.43%  "
$ws.Range("D44").Value = "'105.25"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -6.72%  "
$ws.Range("E45").Value = "  +2.62%  "
$ws.Range("D46").Value = "'0.995"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("E48").Value = "  +3.85%  "
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("E51").Value = "  -3.07%  "
